# Update required/optional fields in cohort template
# - All parent fields marked Optional
# - Mandatory changed to Required
# - Postcode field marked Required
#
# Note: several of the original description strings contain a
# non-breaking space (U+00A0) embedded mid-sentence (an artifact of the
# source template). We preserve that character's position while only
# swapping the specific words the commit touches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# CHILD_SCHOOL_URN: Mandatory -> Required
$ws.Range("A2").Value = "Required, 6 digits, numeric,$($nbsp)use 888888 for school unknown and 999999 for homeschooled "

# CHILD_FIRST_NAME / CHILD_LAST_NAME: Mandatory -> Required
$ws.Range("B2").Value = "Required: Free text$nbsp"
$ws.Range("C2").Value = "Required: Free text$nbsp"

# CHILD_DATE_OF_BIRTH: Mandatory -> Required
$ws.Range("D2").Value = "Required: DD/MM/YYYY"

# CHILD_POSTCODE: Optional -> Required
$ws.Range("I2").Value = "Required: Valid postcode$($nbsp)format"

# PARENT_1_RELATIONSHIP / PARENT_2_RELATIONSHIP: Must be one of -> Optional, must be one of
$ws.Range("L2").Value = "Optional, must be one of:$($nbsp)Mum, Dad, Guardian"
$ws.Range("P2").Value = "Optional, must be one of:$($nbsp)Mum, Dad, Guardian"

# PARENT_2_NAME: Free text -> Optional: Free text
$ws.Range("O2").Value = "Optional: Free text"

# PARENT_2_PHONE: Phone number -> Optional: Phone number
$ws.Range("Q2").Value = "Optional: Phone number$nbsp"

# PARENT_2_EMAIL: Email address -> Optional: Email address
$ws.Range("R2").Value = "Optional: Email address$nbsp"
